$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "vyrobek" column (B) is empty in every data row - drop it entirely.
# Columns shift left: old C (reg.c) -> B, old D (mnozstvi) -> C.
$ws.Columns.Item(2).Delete()

# Row 4: 2025-08-23, regc=100, mnozstvi=100 (plain/default formatting)
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = 45892
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 100

# Row 5: 2025-08-23, regc=102 (keep wrapped style like col B data cells), mnozstvi=200 (plain)
$ws.Range("A2").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 45892
$ws.Range("B2").Copy($ws.Range("B5"))
$ws.Range("B5").Value = 102
$ws.Range("C5").Value = 200

# Row 6: 2025-08-28, regc=394, mnozstvi=4000 (same layout/style as row 2)
$ws.Range("A2:C2").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 45897

# Row 7: 2025-08-28, regc=394, mnozstvi=8000 (same layout/style as row 3)
$ws.Range("A3:C3").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 45897

# Row 8: 2025-08-28, regc=100, mnozstvi=100 (plain/default formatting)
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 45897
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 100

# Row 9: 2025-08-30, regc=102 (wrapped style), mnozstvi=200 (plain)
$ws.Range("A2").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 45899
$ws.Range("B2").Copy($ws.Range("B9"))
$ws.Range("B9").Value = 102
$ws.Range("C9").Value = 200

$ws.Range("I16").Select()
